# feat: add store and warehouse stock in inventory stock report
#
# Adds two new columns ("EN TIENDA" / "EN ALMACEN") between the existing
# "NOMBRE" and "CAJAS" columns, replacing the old single "STOCK" column,
# and re-centers (horizontal + vertical) the title row, the header row and
# the data columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column before the old "CAJAS" column (col H) -------
# This shifts "CAJAS" from H to I, and leaves a blank column H that will
# become "EN ALMACEN". Column G ("STOCK") is renamed to "EN TIENDA" below.
$ws.Columns.Item(8).Insert()

# --- 2. Rewrite the header text --------------------------------------------
$ws.Range("G3").Value = "EN TIENDA"
$ws.Range("H3").Value = "EN ALMACEN"

# new inserted column inherits the width of the column to its left (G),
# matching Excel's default "insert column" behaviour, without the bestFit
# flag (it's a copied width, not an auto-computed one)
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(7).ColumnWidth()

# --- 3. Grow the title merge from A1:G1 to A1:I1 ---------------------------
$ws.Range("A1:I1").Merge()

# --- 4. Re-apply centered alignment --------------------------------------
# Title row / header row keep their own bold+filled styles, now also
# vertically centered.
$ws.Range("A1:I1").HorizontalAlignment = -4108
$ws.Range("A1:I1").VerticalAlignment = -4108

$ws.Range("A3:I3").HorizontalAlignment = -4108
$ws.Range("A3:I3").VerticalAlignment = -4108

# --- 5. Tidy up the view ---------------------------------------------------
$ws.Range("A1").Select()
$ws.Range("A1:I1").EntireColumn.Select()
